# updated MGCC so that BESS is in grid following mode when grid connected
# and seeks target SOC. First draft of MGCC done.

$wb = $excel.ActiveWorkbook

# --- P_req sheet: lower the (negative) requested power magnitude for the
#     first 37 timesteps (rows 2-38, i.e. timesteps 0-36) from -2,000,000
#     to a uniform -1,400,000 ---
$wsP = $wb.Worksheets.Item("P_req")
$wsP.Range("B2:B38").Value = -1400000

# --- try_island sheet: extend the "island mode" flag (value 1) back to
#     cover timesteps 17-36 (rows 19-38), which previously were still 0 ---
$wsIsland = $wb.Worksheets.Item("try_island")
$wsIsland.Range("B19:B38").Value = 1

# --- Update the active sheet / selections to match the author's final
#     view state: try_island becomes the active (selected) tab, with
#     B14:B18 selected; P_req keeps B2:B38 selected as the last range
#     touched on that sheet. ---
$wsP.Range("B2:B38").Select()

$wsIsland.Activate()
$wsIsland.Range("B14:B18").Select()
